$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the top; existing rows 1-12 shift down to 2-13.
$ws.Rows("1").Insert()

# New row 1 is the blank spacer row (single space in every used cell).
$ws.Range("A1:M1").Value = " "

# New row 2 is the "8시" header row (was row 1 before the insert), re-entered
# with a few cells now holding extra blank spaces instead of a single space.
$ws.Range("A2").Value = "8시"
$ws.Range("B2").Value = "   "
$ws.Range("C2").Value = "   "
$ws.Range("D2").Value = "전"
$ws.Range("E2").Value = "차"
$ws.Range("F2").Value = "량"
$ws.Range("G2").Value = "    "
$ws.Range("H2").Value = "수"
$ws.Range("I2").Value = "시"
$ws.Range("J2").Value = "운"
$ws.Range("K2").Value = "행"
$ws.Range("L2").Value = "   "
$ws.Range("M2").Value = " "

# Match the author's final selection.
$ws.Range("M1").Select()
